# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2 through 120) from 2023-10-09 (serial 45208) to 2023-10-13
# (serial 45212), matching the automatic update recorded in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 120
$col = 3  # column C

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $col).Value = 45212
}
